# Updated symbol list on Sun Dec 18 18:50:57 UTC 2022 with GitHub Actions
# Refreshes crypto Price (col D) and Volume(1h) label (col E) cells with
# the latest scraped values, keeping them as plain text (matching the
# original inlineStr cells) and leaving cell styling untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $NewValue) {
    $cell = $Sheet.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

Set-TextValue $ws "D4" "5.523"
Set-TextValue $ws "D5" "0.05616"
Set-TextValue $ws "D7" "0.8079"
Set-TextValue $ws "D8" "1.038"
Set-TextValue $ws "D9" "0.1430"
Set-TextValue $ws "D10" "0.07306"
Set-TextValue $ws "D12" "0.02918"
Set-TextValue $ws "D13" "0.09267"
Set-TextValue $ws "D14" "0.001667"
Set-TextValue $ws "D15" "3.222"
Set-TextValue $ws "D16" "0.04732"
Set-TextValue $ws "D17" "0.0005814"
Set-TextValue $ws "E17" "16OneONE"
Set-TextValue $ws "D18" "0.006388"
Set-TextValue $ws "D19" "0.005070"
Set-TextValue $ws "D20" "0.001058"
Set-TextValue $ws "D21" "0.0001501"
Set-TextValue $ws "D22" "3.986"
Set-TextValue $ws "D23" "3.377"
Set-TextValue $ws "D24" "2.246"
Set-TextValue $ws "D25" "0.3268"
Set-TextValue $ws "D26" "0.1255"
Set-TextValue $ws "E26" "25ProBitTokenPROB"
Set-TextValue $ws "D27" "0.0003302"
Set-TextValue $ws "D40" "0.04145"
Set-TextValue $ws "D41" "0.007013"
Set-TextValue $ws "D42" "0.003502"
Set-TextValue $ws "E42" "41CEJICEJIBestin24h"
Set-TextValue $ws "D43" "0.1040"
Set-TextValue $ws "D44" "0.008570"
Set-TextValue $ws "D45" "0.00005648"
Set-TextValue $ws "D46" "0.00000000750"
Set-TextValue $ws "D47" "0.6804"
Set-TextValue $ws "D48" "0.01620"
Set-TextValue $ws "E48" "47BOLOBOLOWorstin24h"
Set-TextValue $ws "D49" "0.00002101"

Write-Output "Applied 36 cell updates"
